$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.362096
$ws.Range("H2").Value = 1.086288
$ws.Range("I2").Value = 0.01048575317613816
$ws.Range("J2").Value = 0.01048575317613816
$ws.Range("M2").Value = 0.667106
$ws.Range("N2").Value = 2.001318
$ws.Range("O2").Value = 0.003817114239487378
$ws.Range("P2").Value = 0.003817114239487378
$ws.Range("Q2").Value = 0.241556414176
$ws.Range("R2").Value = 2.174007727584
$ws.Range("S2").Value = 0.00004002531776038697
$ws.Range("T2").Value = 0.00004002531776038696

# Row 3
$ws.Range("G3").Value = 0.362096
$ws.Range("H3").Value = 1.086288
$ws.Range("I3").Value = 0.01048575317613816
$ws.Range("J3").Value = 0.01048575317613816
$ws.Range("O3").Value = 0.9945745510447523
$ws.Range("P3").Value = 0.9945745510447522
$ws.Range("Q3").Value = 62.93913336304532
$ws.Range("R3").Value = 566.4522002674079
$ws.Range("S3").Value = 0.0104288632575237
$ws.Range("T3").Value = 0.01042886325752369

# Row 4
$ws.Range("G4").Value = 0.362096
$ws.Range("H4").Value = 1.086288
$ws.Range("I4").Value = 0.01048575317613816
$ws.Range("J4").Value = 0.01048575317613816
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.260372
$ws.Range("N4").Value = 0.7811159999999999
$ws.Range("O4").Value = 0.001489822709979835
$ws.Range("P4").Value = 0.001489822709979834
$ws.Range("Q4").Value = 0.09427965971199999
$ws.Range("R4").Value = 0.8485169374079998
$ws.Range("S4").Value = 0.00001562191321305381
$ws.Range("T4").Value = 0.00001562191321305381

# Row 5
$ws.Range("G5").Value = 0.362096
$ws.Range("H5").Value = 1.086288
$ws.Range("I5").Value = 0.01048575317613816
$ws.Range("J5").Value = 0.01048575317613816
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.020712
$ws.Range("N5").Value = 0.062136
$ws.Range("O5").Value = 0.0001185120057805845
$ws.Range("P5").Value = 0.0001185120057805844
$ws.Range("Q5").Value = 0.007499732351999999
$ws.Range("R5").Value = 0.067497591168
$ws.Range("S5").Value = 0.000001242687641024267
$ws.Range("T5").Value = 0.000001242687641024267

# Row 6
$ws.Range("I6").Value = 0.09151676111574511
$ws.Range("J6").Value = 0.09151676111574511
$ws.Range("M6").Value = 0.667106
$ws.Range("N6").Value = 2.001318
$ws.Range("O6").Value = 0.003817114239487378
$ws.Range("P6").Value = 0.003817114239487378
$ws.Range("Q6").Value = 2.108237747044
$ws.Range("R6").Value = 18.974139723396
$ws.Range("S6").Value = 0.0003493299320066755
$ws.Range("T6").Value = 0.0003493299320066754

# Row 7
$ws.Range("I7").Value = 0.09151676111574511
$ws.Range("J7").Value = 0.09151676111574511
$ws.Range("O7").Value = 0.9945745510447523
$ws.Range("P7").Value = 0.9945745510447522
$ws.Range("S7").Value = 0.09102024159976205
$ws.Range("T7").Value = 0.09102024159976203

# Row 8
$ws.Range("I8").Value = 0.09151676111574511
$ws.Range("J8").Value = 0.09151676111574511
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.260372
$ws.Range("N8").Value = 0.7811159999999999
$ws.Range("O8").Value = 0.001489822709979835
$ws.Range("P8").Value = 0.001489822709979834
$ws.Range("Q8").Value = 0.8228468619279999
$ws.Range("R8").Value = 7.405621757351999
$ws.Range("S8").Value = 0.0001363437490540365
$ws.Range("T8").Value = 0.0001363437490540365

# Row 9
$ws.Range("I9").Value = 0.09151676111574511
$ws.Range("J9").Value = 0.09151676111574511
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.020712
$ws.Range("N9").Value = 0.062136
$ws.Range("O9").Value = 0.0001185120057805845
$ws.Range("P9").Value = 0.0001185120057805844
$ws.Range("Q9").Value = 0.06545559508799999
$ws.Range("R9").Value = 0.5891003557919999
$ws.Range("S9").Value = 0.00001084583492236955
$ws.Range("T9").Value = 0.00001084583492236955

# Row 10
$ws.Range("G10").Value = 21.49951033333333
$ws.Range("H10").Value = 64.498531
$ws.Range("I10").Value = 0.6225933419953967
$ws.Range("J10").Value = 0.6225933419953966
$ws.Range("M10").Value = 0.667106
$ws.Range("N10").Value = 2.001318
$ws.Range("O10").Value = 0.003817114239487378
$ws.Range("P10").Value = 0.003817114239487378
$ws.Range("Q10").Value = 14.34245234042867
$ws.Range("R10").Value = 129.082071063858
$ws.Range("S10").Value = 0.002376509911140664
$ws.Range("T10").Value = 0.002376509911140663

# Row 11
$ws.Range("G11").Value = 21.49951033333333
$ws.Range("H11").Value = 64.498531
$ws.Range("I11").Value = 0.6225933419953967
$ws.Range("J11").Value = 0.6225933419953966
$ws.Range("O11").Value = 0.9945745510447523
$ws.Range("P11").Value = 0.9945745510447522
$ws.Range("Q11").Value = 3737.021530505274
$ws.Range("R11").Value = 33633.19377454747
$ws.Range("S11").Value = 0.6192154935985237
$ws.Range("T11").Value = 0.6192154935985235

# Row 12
$ws.Range("G12").Value = 21.49951033333333
$ws.Range("H12").Value = 64.498531
$ws.Range("I12").Value = 0.6225933419953967
$ws.Range("J12").Value = 0.6225933419953966
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.260372
$ws.Range("N12").Value = 0.7811159999999999
$ws.Range("O12").Value = 0.001489822709979835
$ws.Range("P12").Value = 0.001489822709979834
$ws.Range("Q12").Value = 5.597870504510666
$ws.Range("R12").Value = 50.380834540596
$ws.Range("S12").Value = 0.000927553699986984
$ws.Range("T12").Value = 0.0009275536999869836

# Row 13
$ws.Range("G13").Value = 21.49951033333333
$ws.Range("H13").Value = 64.498531
$ws.Range("I13").Value = 0.6225933419953967
$ws.Range("J13").Value = 0.6225933419953966
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.020712
$ws.Range("N13").Value = 0.062136
$ws.Range("O13").Value = 0.0001185120057805845
$ws.Range("P13").Value = 0.0001185120057805844
$ws.Range("Q13").Value = 0.445297858024
$ws.Range("R13").Value = 4.007680722216
$ws.Range("S13").Value = 0.00007378478574551185
$ws.Range("T13").Value = 0.00007378478574551184

# Row 14
$ws.Range("G14").Value = 9.510307666666666
$ws.Range("H14").Value = 28.530923
$ws.Range("I14").Value = 0.27540414371272
$ws.Range("J14").Value = 0.27540414371272
$ws.Range("M14").Value = 0.667106
$ws.Range("N14").Value = 2.001318
$ws.Range("O14").Value = 0.003817114239487378
$ws.Range("P14").Value = 0.003817114239487378
$ws.Range("Q14").Value = 6.344383306279332
$ws.Range("R14").Value = 57.09944975651399
$ws.Range("S14").Value = 0.001051249078579652
$ws.Range("T14").Value = 0.001051249078579652

# Row 15
$ws.Range("G15").Value = 9.510307666666666
$ws.Range("H15").Value = 28.530923
$ws.Range("I15").Value = 0.27540414371272
$ws.Range("J15").Value = 0.27540414371272
$ws.Range("O15").Value = 0.9945745510447523
$ws.Range("P15").Value = 0.9945745510447522
$ws.Range("Q15").Value = 1653.071347255771
$ws.Range("R15").Value = 14877.64212530194
$ws.Range("S15").Value = 0.273909952588943
$ws.Range("T15").Value = 0.273909952588943

# Row 16
$ws.Range("G16").Value = 9.510307666666666
$ws.Range("H16").Value = 28.530923
$ws.Range("I16").Value = 0.27540414371272
$ws.Range("J16").Value = 0.27540414371272
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.260372
$ws.Range("N16").Value = 0.7811159999999999
$ws.Range("O16").Value = 0.001489822709979835
$ws.Range("P16").Value = 0.001489822709979834
$ws.Range("Q16").Value = 2.476217827785333
$ws.Range("R16").Value = 22.285960450068
$ws.Range("S16").Value = 0.0004103033477257604
$ws.Range("T16").Value = 0.0004103033477257603

# Row 17
$ws.Range("G17").Value = 9.510307666666666
$ws.Range("H17").Value = 28.530923
$ws.Range("I17").Value = 0.27540414371272
$ws.Range("J17").Value = 0.27540414371272
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.020712
$ws.Range("N17").Value = 0.062136
$ws.Range("O17").Value = 0.0001185120057805845
$ws.Range("P17").Value = 0.0001185120057805844
$ws.Range("Q17").Value = 0.196977492392
$ws.Range("R17").Value = 1.772797431528
$ws.Range("S17").Value = 0.00003263869747167879
$ws.Range("T17").Value = 0.00003263869747167879
